$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.288.69'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.831.39'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6036'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.68%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07056'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.92%  '
$ws.Range('E9').Value = '  -3.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.56'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07657'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.80%  '
$ws.Range('D12').Value = '1.828.82'
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.795'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6292'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009883'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.19%  '
$ws.Range('D16').Value = '2.085.08'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.16'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.855'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.70%  '
$ws.Range('D19').Value = '29.297.80'
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '224.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.004'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.72'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.011'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.85%  '
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1305'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.994'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.486'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06375'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -12.75%  '
$ws.Range('E31').Value = '  -2.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.847'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.804'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.86%  '
$ws.Range('E34').Value = '  -2.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.734'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6477'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.546'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('D38').Value = '1.217.84'
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('E39').Value = '  -2.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01750'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.554'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.12%  '
$ws.Range('E42').Value = '  -6.01%  '
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').Value = '1.999.99'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '62.82'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000117'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.575'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.583'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.67%  '
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05506'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.65%  '
